# Login Authentication Problems Fixed (Team Commit #1)
# Update the Week 8-11 planning rows with corrected date ranges / task
# assignments, and move the sheet's active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 8 row (row 9): new date range, clear the task assigned in col B ---
$ws.Range("A9").Value = "Week 8 (6.07.2020 - 12.07.2020)"
$ws.Range("B9").Value = $null

# --- Week 9 row (row 10): new date range, reassign task ---
$ws.Range("A10").Value = "Week 9 (03.08.2020 - 09.08.2020)"
$ws.Range("B10").Value = "Permissions"

# --- Week 10 row (row 11): new date range, reassign task ---
$ws.Range("A11").Value = "Week 10 (10.08.2020-16.08.2020)"
$ws.Range("B11").Value = "Leave System"

# --- Week 11 row (row 12): new date range (task in col B unchanged) ---
$ws.Range("A12").Value = "Week 11 (17.08.2020 - 23.08.2020)"

# --- Update the view: scroll back to top and move the selection to B11 ---
$ws.Range("B11").Select()
